# MarshalByRef works for return types
#
# Adds a third "Tests for Factory" style scenario (rows 17-19) exercising
# _xll.dnaFactory, driven off three labeled inputs (One/Two/Three -> 1/2/3)
# in columns C:D, with the factory call living in A17 as an array formula
# spanning A17:A19 (mirroring the other _xll.* array-formula blocks already
# on the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New inputs for the factory test (rows 17-19, columns C/D) ----------
$ws.Range("C17").Value = "One"
$ws.Range("D17").Value = 1

$ws.Range("C18").Value = "Two"
$ws.Range("D18").Value = 2

$ws.Range("C19").Value = "Three"
$ws.Range("D19").Value = 3

# Seed A18 with the long diagnostic text so the column-A width picks up its
# length before we AutoFit; the array formula written below becomes the
# permanent contents of A17:A19 (matching the other _xll array blocks).
$ws.Range("A18").Value = "Exception has been thrown by the target of an invocation."
$ws.Columns.Item(1).AutoFit()

$ws.Range("A17:A19").FormulaArray = "=_xll.dnaFactory(C17:C19,D17:D19)"

# --- Selection moves to the newly added row ------------------------------
$ws.Range("A19").Select()
